$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. First paragraph: append "  (This is a change – Version for
#    branch alternate)" as three separate red (C00000) runs after the
#    existing plain-text run (which grows a trailing two spaces).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "This is a Microsoft word document.  "

$afterFirst = $d.Paragraphs(1).Range.End - 1
$enDash = [char]0x2013

$r2 = $d.Range($afterFirst, $afterFirst)
$r2.InsertAfter("(This is a change " + $enDash + " Ve")
$r2.Font.Color = 192

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter("rsion for branch alternate")
$r3.Font.Color = 192

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(")")
$r4.Font.Color = 192

# ------------------------------------------------------------------
# 2. Add a new, empty, shaded paragraph after the final paragraph of
#    the speech ("... we are free at last.").
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute("Thank God almighty, we are free at last.", $true, $false, $false, $false, $false, $true, 1, $false, "Thank God almighty, we are free at last.^p", 2)

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Style = "Normal"
$newPara.Format.Shading.Texture = 0
$newPara.Format.Shading.ForegroundPatternColor = -16777216
$newPara.Format.Shading.BackgroundPatternColor = 0xF9F9F9

# ------------------------------------------------------------------
# 3. Drop the handful of unused styles that this edit session's
#    re-save no longer carries forward (Word prunes styles that
#    nothing in the document references once it is touched).  Delete
#    from the highest original index down so earlier deletes don't
#    shift the indices of styles still pending removal.
# ------------------------------------------------------------------
$unusedStyles = @(
  "podcast-toolssubscribe-links",
  "generic-title",
  "subscribe-more-info",
  "subscribe",
  "audio-tool",
  "Heading4Char",
  "Heading2Char",
  "Hyperlink",
  "apple-converted-space",
  "Heading4",
  "Heading2"
)
foreach ($styleName in $unusedStyles) {
  $d.Styles($styleName).Delete()
}
